$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D23").Value = "안녕하세요! 통-하!`nR에 좀 더 익숙한 상황에서 파이썬으로 시계열분석을 해야해서 참고할 만한 파이썬 시계열 분석 책이 있는 지 여쭤보고자 합"
$ws.Range("E23").Value = "https://theonly1.tistory.com/2707"

$ws.Range("D44").Value = "2020년도 IT Start-up와 산업 동향 파악 및 전망 (1)"
$ws.Range("E44").Value = "https://engineering-ladder.tistory.com/76"
